# Adds three new rows (20-22) of "Child Development" (Young et al., 2022)
# visualization entries to the dataset, matching how the data was pasted
# in column-by-column (A, B, C, D share the same repeated values across
# the three rows; E/F/G/H/I/J/K/L vary per row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$journal = "Child Development"
$inText = "Young et al., 2022"
$title = "Hidden talents in context: Cognitive performance with abstract versus ecological stimuli among adversity-exposed youth."
$link = "/publications/journal/young-2022-cd-1/"

# Column A - journal (same for all three new rows)
$ws.Range("A20").Value = $journal
$ws.Range("A21").Value = $journal
$ws.Range("A22").Value = $journal

# Column B - in_text citation
$ws.Range("B20").Value = $inText
$ws.Range("B21").Value = $inText
$ws.Range("B22").Value = $inText

# Column C - title
$ws.Range("C20").Value = $title
$ws.Range("C21").Value = $title
$ws.Range("C22").Value = $title
# Column C cells carry the default (unstyled) format rather than the
# column-inherited style used elsewhere in the sheet.
$ws.Range("C20:C22").Style = "Normal"

# Column D - link
$ws.Range("D20").Value = $link
$ws.Range("D21").Value = $link
$ws.Range("D22").Value = $link

# Column E - figure number
$ws.Range("E20").Value = 1
$ws.Range("E21").Value = 2
$ws.Range("E22").Value = 3

# Column F - file
$ws.Range("F20").Value = "plots/cd2022-figure1.jpg"
$ws.Range("F21").Value = "plots/cd2022-figure2.jpg"
$ws.Range("F22").Value = "plots/cd2022-figure3.jpg"

# Column G - caption
$ws.Range("G20").Value = "Figure 1. Schematic display of the attention shifting and working memory updating tasks: (a) abstract attention shifting, (b) ecological attention shifting, (c) abstract working memory updating, and (d) ecological working memory updating"
$ws.Range("G21").Value = "Figure 2. Visualization of multiverse attention-shifting results. Unpredictability and violence exposure multiverse analyses were confirmatory and poverty analyses were exploratory. (a) visualizes the multiverse task version × adversity interaction on abstract and ecological task versions across high (+1 SD) and low (−1 SD) adversity exposure (y-axis was reversed so that higher values = faster shifting), (b) plots p-curves associated with each interaction term, (c) plots sorted interaction β-coefficients across each arbitrary decision, (d) plots the sample sizes for each effect, and (e) is a specification grid indicating the data processing decisions associated with each effect"
$ws.Range("G22").Value = "Figure 3. Visualization of multiverse working memory updating results. Unpredictability and Violence exposure multiverse analyses were confirmatory and poverty analyses were exploratory. (a) visualizes the multiverse task version × adversity interaction on abstract and ecological task versions across high (+1 SD) and low (−1 SD) adversity exposure, (b) plots p-curves associated with each interaction term, (c) plots sorted interaction β-coefficients across each arbitrary decision, (d) plots the sample sizes for each effect, and (e) is a specification grid indicating the data processing decisions associated with each effect. Proportions of each arbitrary decision with p-values < .05 are indicated on the right side of each specification grid. Blank proportions indicate proportions = 0. Teal lines and points reflect individual multiverse effect sizes with p-values < .05"

# Column H - tool
$ws.Range("H20").Value = "powerpoint"
$ws.Range("H21").Value = "ggplot2"
$ws.Range("H22").Value = "ggplot2"

# Column I - width
$ws.Range("I20").Value = 2250
$ws.Range("I21").Value = 1950
$ws.Range("I22").Value = 1950

# Column J - height
$ws.Range("J20").Value = 1500
$ws.Range("J21").Value = 2175
$ws.Range("J22").Value = 2175

# Column K - plot_type
$ws.Range("K20").Value = "conceptual"
$ws.Range("K21").Value = "data"
$ws.Range("K22").Value = "data"

# Column L - featured
$ws.Range("L20").Value = 1
$ws.Range("L21").Value = 1
$ws.Range("L22").Value = 1

# Match the resulting selection left on the sheet
$ws.Range("J24").Select() | Out-Null
